$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price-snapshot column is inserted right before the existing
# "nom" column (AF), pushing "nom" -> AG and "url_produit" -> AH. This
# mirrors the other timestamped snapshot columns (B..AE) that already
# exist for every prior scrape run.
$ws.Columns("AF").Insert()

# Header for the newly inserted column (AF, column 32) is the new scrape
# timestamp.
$ws.Cells.Item(1, 32).Value = "2026-01-29 03:25:05"

# For every product row, carry the most recently recorded price (the
# previous last snapshot column, now AE / column 31) forward into the
# new column, same as every other repeated-price snapshot column in the
# sheet. Rows that never had a recorded price (AE blank) stay blank in
# the new column too.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $priorPrice = $ws.Cells.Item($r, 31).Value()
    if ($priorPrice -is [double]) {
        $ws.Cells.Item($r, 32).Value = $priorPrice
    }
}
